$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: replace name "vivek shree unikrishan" -> "Johnson powder boi"
# and convert the id column from text "EYEQ202504215" to a plain number 202504215
$ws.Range("A2").Value = "Johnson powder boi"
$ws.Range("B2").Value = 202504215

# Row 3: id EYEQ202504216 -> numeric 202504216
$ws.Range("B3").Value = 202504216

# Row 4: id EYEQ202504217 -> numeric 202504217
$ws.Range("B4").Value = 202504217

# Row 5: id EYEQ202504218 -> numeric 202504218
$ws.Range("B5").Value = 202504218

# Update the active selection to the whole A1:A5 column, active cell A1
$ws.Range("A1:A5").Select()
